# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / handoff / handback
# datetime-stamp cells that get refreshed whenever the handback report
# is regenerated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the cad4860f row.
$wsOverview.Range("G4").Value = "2016-08-22 00:56:33"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for the cad4860f row.
$wsZhCn.Range("H4").Value = "2016-08-22 00:56:28"
$wsZhCn.Range("K4").Value = "2016-08-22 00:56:45"

# de-de sheet: Correspond Handoff Datetime (mirrors the Overview value) /
# Correspond Handback DateTime for the cad4860f row.
$wsDeDe.Range("H4").Value = "2016-08-22 00:56:33"
$wsDeDe.Range("K4").Value = "2016-08-22 00:56:51"
